$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestPlan")

# Swap the "Execute" flag (column D) between row 2 and row 3
$ws.Range("D2").Value = "Y"
$ws.Range("D3").Value = "N"

# Update the "Rounding Precision" value (column I) on row 3
$ws.Range("I3").Value = "0.001"

# Update the active selection to reflect the last edited cell
$ws.Range("I3").Select()
